$wb = $excel.ActiveWorkbook

# Grab references to the sheets we need to touch
$wsParameters   = $wb.Worksheets.Item("Parameters")
$wsCompartments = $wb.Worksheets.Item("Compartments")
$wsBulkCargo    = $wb.Worksheets.Item("BulkCargo")
$wsGeneralCargo = $wb.Worksheets.Item("GeneralCargo")

# --- Compartments: change several "max" dropdown selections to "current" ---
$wsCompartments.Range("E24:E27").Value = "current"
$wsCompartments.Range("E32:E37").Value = "current"
$wsCompartments.Range("E39").Value = "current"
$wsCompartments.Range("E42:E44").Value = "current"
$wsCompartments.Range("E48:E52").Value = "current"

# --- GeneralCargo: clear the placeholder "-" values on row 2 (A2:K2), keep styles ---
$wsGeneralCargo.Range("A2:K2").ClearContents()

# --- Update each sheet's selection (activating each sheet briefly updates its
#     own stored selection/view state), finishing on Compartments so that it
#     ends up as the workbook's active sheet/tab. ---
$wsBulkCargo.Activate()
$wsBulkCargo.Range("H22").Select()

$wsGeneralCargo.Activate()
$wsGeneralCargo.Range("K14").Select()

$wsCompartments.Activate()
$wsCompartments.Range("I48").Select()

Write-Output "done"
